$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.815046728013995
$ws.Range("C2").Value = 5.70633596196627
$ws.Range("D2").Value = 4.388557147948176
$ws.Range("F2").Value = 20.57419389904199
$ws.Range("G2").Value = 22.33688830220235
$ws.Range("H2").Value = 12.72332909067796
$ws.Range("I2").Value = 18.26488392962155
$ws.Range("K2").Value = 8.207227890698872
$ws.Range("M2").Value = 19.9702627975824
$ws.Range("N2").Value = 17.52886423546707
$ws.Range("O2").Value = 18.49877040484611
$ws.Range("B3").Value = 6.497506685765712
$ws.Range("C3").Value = 5.639841806035915
$ws.Range("D3").Value = 4.292823474867044
$ws.Range("F3").Value = 20.58850134445391
$ws.Range("G3").Value = 22.35810170889689
$ws.Range("H3").Value = 12.75971040644842
$ws.Range("I3").Value = 18.334547005731
$ws.Range("K3").Value = 8.02693027849134
$ws.Range("M3").Value = 19.3686812694531
$ws.Range("N3").Value = 17.58741963363986
$ws.Range("O3").Value = 18.55277244908273
$ws.Range("B4").Value = 6.292897430294553
$ws.Range("C4").Value = 5.598179885013691
$ws.Range("D4").Value = 4.232111532763306
$ws.Range("F4").Value = 20.60247818255181
$ws.Range("G4").Value = 22.37825398275167
$ws.Range("H4").Value = 12.78379279349314
$ws.Range("I4").Value = 18.38049201089119
$ws.Range("K4").Value = 7.912804917286582
$ws.Range("M4").Value = 18.9977209924962
$ws.Range("N4").Value = 17.62494880326924
$ws.Range("O4").Value = 18.58951484654808
$ws.Range("B5").Value = 6.207155235874332
$ws.Range("C5").Value = 5.58100310963731
$ws.Range("D5").Value = 4.206903565010457
$ws.Range("F5").Value = 20.60947859164298
$ws.Range("G5").Value = 22.38825359897006
$ws.Range("H5").Value = 12.79404535491305
$ws.Range("I5").Value = 18.40001218088663
$ws.Range("K5").Value = 7.865475436861673
$ws.Range("M5").Value = 18.84640825553386
$ws.Range("N5").Value = 17.64064010491474
$ws.Range("O5").Value = 18.60538781691321
$ws.Range("B6").Value = 6.192777138102747
$ws.Range("C6").Value = 5.578139200691145
$ws.Range("D6").Value = 4.202690128599958
$ws.Range("F6").Value = 20.61071976949599
$ws.Range("G6").Value = 22.39002182045977
$ws.Range("H6").Value = 12.79577429349203
$ws.Range("I6").Value = 18.40330163068156
$ws.Range("K6").Value = 7.857567892835828
$ws.Range("M6").Value = 18.82128194377804
$ws.Range("N6").Value = 17.64326970970993
$ws.Range("O6").Value = 18.60807783843991
$ws.Range("B7").Value = 6.291750556965929
$ws.Range("C7").Value = 5.597949025494041
$ws.Range("D7").Value = 4.231773436720113
$ws.Range("F7").Value = 20.60256731132868
$ws.Range("G7").Value = 22.37838161116451
$ws.Range("H7").Value = 12.78392928618558
$ws.Range("I7").Value = 18.38075203930881
$ws.Range("K7").Value = 7.912169893666126
$ws.Range("M7").Value = 18.99568056979439
$ws.Range("N7").Value = 17.62515880868274
$ws.Range("O7").Value = 18.58972527213907
$ws.Range("B8").Value = 6.707593030972865
$ws.Range("C8").Value = 5.683588196532627
$ws.Range("D8").Value = 4.355959554726165
$ws.Range("F8").Value = 20.57804919311189
$ws.Range("G8").Value = 22.34272090556635
$ws.Range("H8").Value = 12.73551151245171
$ws.Range("I8").Value = 18.28824530297574
$ws.Range("K8").Value = 8.145793564847722
$ws.Range("M8").Value = 19.76331608163918
$ws.Range("N8").Value = 17.5487280139388
$ws.Range("O8").Value = 18.51664548844293
$ws.Range("B9").Value = 7.444381231060386
$ws.Range("C9").Value = 5.844455419620427
$ws.Range("D9").Value = 4.583424229520397
$ws.Range("F9").Value = 20.5711761033821
$ws.Range("G9").Value = 22.32950421458279
$ws.Range("H9").Value = 12.65439316384579
$ws.Range("I9").Value = 18.13201917438983
$ws.Range("K9").Value = 8.575228381192419
$ws.Range("M9").Value = 21.24554056551027
$ws.Range("N9").Value = 17.41127810301824
$ws.Range("O9").Value = 18.40183771785104
$ws.Range("B10").Value = 7.935377701262437
$ws.Range("C10").Value = 5.957748895119175
$ws.Range("D10").Value = 4.739839182718144
$ws.Range("F10").Value = 20.59122213115236
$ws.Range("G10").Value = 22.35451585171891
$ws.Range("H10").Value = 12.60321210214591
$ws.Range("I10").Value = 18.03260915605867
$ws.Range("K10").Value = 8.871331953332106
$ws.Range("M10").Value = 22.30749231459095
$ws.Range("N10").Value = 17.31776677671672
$ws.Range("O10").Value = 18.33493815727597
$ws.Range("B11").Value = 8.147459601278678
$ws.Range("C11").Value = 6.00810769963093
$ws.Range("D11").Value = 4.808502485489333
$ws.Range("F11").Value = 20.60577076146891
$ws.Range("G11").Value = 22.37343934784858
$ws.Range("H11").Value = 12.58175321156163
$ws.Range("I11").Value = 17.99072682809769
$ws.Range("K11").Value = 9.001468053454982
$ws.Range("M11").Value = 22.78233578281947
$ws.Range("N11").Value = 17.27682609354684
$ws.Range("O11").Value = 18.30830735837555
$ws.Range("B12").Value = 8.226125433879119
$ws.Range("C12").Value = 6.026998205906953
$ws.Range("D12").Value = 4.834132582110074
$ws.Range("F12").Value = 20.61205773363071
$ws.Range("G12").Value = 22.38168847531701
$ws.Range("H12").Value = 12.57388932554865
$ws.Range("I12").Value = 17.97534782378479
$ws.Range("K12").Value = 9.05006432429426
$ws.Range("M12").Value = 22.96078560850403
$ws.Range("N12").Value = 17.26155099791131
$ws.Range("O12").Value = 18.29877077804568
$ws.Range("B13").Value = 8.209256828222941
$ws.Range("C13").Value = 6.022937922318159
$ws.Range("D13").Value = 4.828629388123741
$ws.Range("F13").Value = 20.61066918886039
$ws.Range("G13").Value = 22.37986374544305
$ws.Range("H13").Value = 12.57557129748227
$ws.Range("I13").Value = 17.97863856748709
$ws.Range("K13").Value = 9.039629055305912
$ws.Range("M13").Value = 22.92241651172361
$ws.Range("N13").Value = 17.26483063618093
$ws.Range("O13").Value = 18.30080026430717
$ws.Range("B14").Value = 8.15396457367827
$ws.Range("C14").Value = 6.009665480559086
$ws.Range("D14").Value = 4.810618599431809
$ws.Range("F14").Value = 20.60627242789408
$ws.Range("G14").Value = 22.37409631463606
$ws.Range("H14").Value = 12.58110099205699
$ws.Range("I14").Value = 17.98945194584143
$ws.Range("K14").Value = 9.005479943631043
$ws.Range("M14").Value = 22.79704523091957
$ws.Range("N14").Value = 17.2755648366111
$ws.Range("O14").Value = 18.30751179050232
$ws.Range("B15").Value = 8.119881623120257
$ws.Range("C15").Value = 6.001512095966385
$ws.Range("D15").Value = 4.799537760299351
$ws.Range("F15").Value = 20.60368046382208
$ws.Range("G15").Value = 22.37070458546194
$ws.Range("H15").Value = 12.58452222358148
$ws.Range("I15").Value = 17.99613810292411
$ws.Range("K15").Value = 8.984472849702636
$ws.Range("M15").Value = 22.7200692142139
$ws.Range("N15").Value = 17.28216952345553
$ws.Range("O15").Value = 18.31169418359076
$ws.Range("B16").Value = 7.921288639885542
$ws.Range("C16").Value = 5.954433263851292
$ws.Range("D16").Value = 4.735300615628509
$ws.Range("F16").Value = 20.59038035800446
$ws.Range("G16").Value = 22.35343086664999
$ws.Range("H16").Value = 12.60465117819586
$ws.Range("I16").Value = 18.03541352638633
$ws.Range("K16").Value = 8.862733050284003
$ws.Range("M16").Value = 22.27627870391644
$ws.Range("N16").Value = 17.32047437178367
$ws.Range("O16").Value = 18.33675515460866
$ws.Range("B17").Value = 7.796552146661381
$ws.Range("C17").Value = 5.925242918482562
$ws.Range("D17").Value = 4.695245661195078
$ws.Range("F17").Value = 20.58360989454716
$ws.Range("G17").Value = 22.34476572046909
$ws.Range("H17").Value = 12.61746664033268
$ws.Range("I17").Value = 18.06036359548493
$ws.Range("K17").Value = 8.786860812628836
$ws.Range("M17").Value = 22.001784180187
$ws.Range("N17").Value = 17.34438136186636
$ws.Range("O17").Value = 18.3531039391557
$ws.Range("B18").Value = 7.723746858487164
$ws.Range("C18").Value = 5.908342995178542
$ws.Range("D18").Value = 4.671973389993168
$ws.Range("F18").Value = 20.58022699238306
$ws.Range("G18").Value = 22.34049217746954
$ws.Range("H18").Value = 12.62500941037468
$ws.Range("I18").Value = 18.07502853113309
$ws.Range("K18").Value = 8.742793113479406
$ws.Range("M18").Value = 21.84313574337649
$ws.Range("N18").Value = 17.35828254279016
$ws.Range("O18").Value = 18.36286513886711
$ws.Range("B19").Value = 7.698914955424161
$ws.Range("C19").Value = 5.902602311426532
$ws.Range("D19").Value = 4.66405407260381
$ws.Range("F19").Value = 20.5791694946766
$ws.Range("G19").Value = 22.33916727853104
$ws.Range("H19").Value = 12.62759274784559
$ws.Range("I19").Value = 18.08004779031709
$ws.Range("K19").Value = 8.727799908197197
$ws.Range("M19").Value = 21.7892941729107
$ws.Range("N19").Value = 17.36301514163209
$ws.Range("O19").Value = 18.36623153024982
$ws.Range("B20").Value = 7.809940467821311
$ws.Range("C20").Value = 5.928361779175539
$ws.Range("D20").Value = 4.699533866682298
$ws.Range("F20").Value = 20.58427772511873
$ws.Range("G20").Value = 22.34561462338402
$ws.Range("H20").Value = 12.61608464762225
$ws.Range("I20").Value = 18.05767508277654
$ws.Range("K20").Value = 8.794982042618942
$ws.Range("M20").Value = 22.03108525801113
$ws.Range("N20").Value = 17.34182085571151
$ws.Range("O20").Value = 18.35132654218619
$ws.Range("B21").Value = 8.170250059147595
$ws.Range("C21").Value = 6.013568859955295
$ws.Range("D21").Value = 4.81591897777867
$ws.Range("F21").Value = 20.60754278405679
$ws.Range("G21").Value = 22.37576097208438
$ws.Range("H21").Value = 12.57946967342684
$ws.Range("I21").Value = 17.98626273695202
$ws.Range("K21").Value = 9.01552912055738
$ws.Range("M21").Value = 22.83390813780775
$ws.Range("N21").Value = 17.27240575940581
$ws.Range("O21").Value = 18.30552557308283
$ws.Range("B22").Value = 8.396137145121415
$ws.Range("C22").Value = 6.068207363274011
$ws.Range("D22").Value = 4.889814399258841
$ws.Range("F22").Value = 20.62727905507231
$ws.Range("G22").Value = 22.40177533877164
$ws.Range("H22").Value = 12.55706751663458
$ws.Range("I22").Value = 17.94239415669161
$ws.Range("K22").Value = 9.155674635604502
$ws.Range("M22").Value = 23.35058404071782
$ws.Range("N22").Value = 17.22836868970705
$ws.Range("O22").Value = 18.27878620099404
$ws.Range("B23").Value = 8.276461651825253
$ws.Range("C23").Value = 6.039144899836333
$ws.Range("D23").Value = 4.850577520022016
$ws.Range("F23").Value = 20.61633203151037
$ws.Range("G23").Value = 22.38731438519773
$ws.Range("H23").Value = 12.5688842064185
$ws.Range("I23").Value = 17.96555089703722
$ws.Range("K23").Value = 9.081250302938166
$ws.Range("M23").Value = 23.07561196366596
$ws.Range("N23").Value = 17.25175094731483
$ws.Range("O23").Value = 18.29276487078497
$ws.Range("B24").Value = 7.803891009238301
$ws.Range("C24").Value = 5.926952109072104
$ws.Range("D24").Value = 4.697595928434599
$ws.Range("F24").Value = 20.58397421163099
$ws.Range("G24").Value = 22.345228628477
$ws.Range("H24").Value = 12.6167089014635
$ws.Range("I24").Value = 18.05888956017525
$ws.Range("K24").Value = 8.791311827051986
$ws.Range("M24").Value = 22.01784083686032
$ws.Range("N24").Value = 17.34297797223591
$ws.Range("O24").Value = 18.35212897569277
$ws.Range("B25").Value = 7.253755656888031
$ws.Range("C25").Value = 5.801755054863947
$ws.Range("D25").Value = 4.523707206124109
$ws.Range("F25").Value = 20.56862371861476
$ws.Range("G25").Value = 22.32698788862977
$ws.Range("H25").Value = 12.6748587925013
$ws.Range("I25").Value = 18.17158493683058
$ws.Range("K25").Value = 8.462333409599259
$ws.Range("M25").Value = 20.84843605653291
$ws.Range("N25").Value = 17.44714200147306
$ws.Range("O25").Value = 18.42983699214743
